$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WebViewMeasurement")
$ws.Range("B10:C19").ClearContents()
$ws.Range("E10:F19").ClearContents()
$ws.Range("H10:I19").ClearContents()
$wb.RefreshAll()
